$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The glossary had duplicate "Frequency"/"Severity" terms (capitalised) that
# duplicated the lower-case "frequency"/"severity" terms used elsewhere in
# the sheet. Clean up by normalising the two glossary entries to lower case,
# matching the rest of the table.
$ws.Range("A3").Value = "frequency"
$ws.Range("A4").Value = "severity"

# Leave the cursor where the author finished editing.
$ws.Range("B7").Select()
